$d = $word.ActiveDocument

$replacements = @(
    @{old = "71×41=2911"; new = "49×20=980"},
    @{old = "36×74=2664"; new = "80×23=1840"},
    @{old = "74×42=3108"; new = "77×40=3080"},
    @{old = "59×79=4661"; new = "93×14=1302"},
    @{old = "91×72=6552"; new = "49×56=2744"},
    @{old = "31×70=2170"; new = "87×91=7917"},
    @{old = "96×33=3168"; new = "61×25=1525"},
    @{old = "57×92=5244"; new = "20×90=1800"},
    @{old = "44×93=4092"; new = "48×41=1968"},
    @{old = "25×13=325"; new = "62×36=2232"},
    @{old = "65×89=5785"; new = "41×37=1517"},
    @{old = "88×47=4136"; new = "71×62=4402"},
    @{old = "86×85=7310"; new = "37×17=629"},
    @{old = "89×17=1513"; new = "67×97=6499"},
    @{old = "44×52=2288"; new = "29×22=638"},
    @{old = "57×75=4275"; new = "83×41=3403"},
    @{old = "39×52=2028"; new = "52×70=3640"},
    @{old = "36×48=1728"; new = "13×92=1196"},
    @{old = "19×44=836"; new = "21×63=1323"},
    @{old = "94×86=8084"; new = "58×78=4524"},
    @{old = "99×94=9306"; new = "71×33=2343"},
    @{old = "92×40=3680"; new = "74×51=3774"},
    @{old = "56×23=1288"; new = "64×84=5376"},
    @{old = "24×90=2160"; new = "56×22=1232"},
    @{old = "33×11=363"; new = "77×60=4620"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
